$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new "actual hours" data for row 20 (columns H:L)
$ws.Range("H20:L20").Value = 4

# Highlight the newly entered cells with a yellow fill (BGR Long for pure yellow)
$ws.Range("H20:L20").Interior.Color = 65535

# Move the selection to where the user ended up
$ws.Range("O24").Select()
